$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (incl. date number format / style) from the row above
# so the new date cell keeps the same style index as existing date cells.
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)

$ws.Range("A5").Value = 41426
$ws.Range("B5").Value = 0.75
$ws.Range("D5").Value = "Implementation design of new sync objects "

$ws.Range("D5").Select()
